# Weekly fruit/vegetable price update: insert two new data rows
# (a new week's worth of price observations) into the "Acelga" sheet,
# right before the row that used to be row 256.
# This pushes the former rows 256..361 down to 258..363.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 256-257 (EntireRow.Insert shifts rows 256..361 down to 258..363)
$ws.Range("A256:A257").EntireRow.Insert()

# New row 256
$ws.Cells.Item(256, 1).Value = 10
$ws.Cells.Item(256, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(256, 3).Value = 'La Araucanía'
$ws.Cells.Item(256, 4).Value = 44784
$ws.Cells.Item(256, 5).Value = 9
$ws.Cells.Item(256, 6).Value = 100112009
$ws.Cells.Item(256, 7).Value = 'Acelga'
$ws.Cells.Item(256, 8).Value = 'Sin especificar'
$ws.Cells.Item(256, 9).Value = 'Primera'
$ws.Cells.Item(256, 10).Value = 80
$ws.Cells.Item(256, 11).Value = 10000
$ws.Cells.Item(256, 12).Value = 10000
$ws.Cells.Item(256, 13).Value = 10000
$ws.Cells.Item(256, 14).Value = '$/docena de atados (12 kilos)'
$ws.Cells.Item(256, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(256, 16).Value = 833
$ws.Cells.Item(256, 17).Value = 12
$ws.Cells.Item(256, 18).Value = 'Hortaliza'

# New row 257
$ws.Cells.Item(257, 1).Value = 10
$ws.Cells.Item(257, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(257, 3).Value = 'La Araucanía'
$ws.Cells.Item(257, 4).Value = 44784
$ws.Cells.Item(257, 5).Value = 9
$ws.Cells.Item(257, 6).Value = 100112009
$ws.Cells.Item(257, 7).Value = 'Acelga'
$ws.Cells.Item(257, 8).Value = 'Sin especificar'
$ws.Cells.Item(257, 9).Value = 'Primera'
$ws.Cells.Item(257, 10).Value = 40
$ws.Cells.Item(257, 11).Value = 6000
$ws.Cells.Item(257, 12).Value = 6000
$ws.Cells.Item(257, 13).Value = 6000
$ws.Cells.Item(257, 14).Value = '$/docena de atados (6 kilos)'
$ws.Cells.Item(257, 15).Value = 'Región del Maule'
$ws.Cells.Item(257, 16).Value = 1000
$ws.Cells.Item(257, 17).Value = 6
$ws.Cells.Item(257, 18).Value = 'Hortaliza'
